$d = $word.ActiveDocument

$bodyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>Circle Language Spec Plan</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>Orient in First Four Principles</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>2008-03</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Project </w:t>
      </w:r>
      <w:r>
        <w:t>Summary</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="142"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="142"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">Author: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>JJ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> van Zon</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="142"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>Location: Oosterhout</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="20"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>, The Netherlands</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="142"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
      </w:pPr>
      <w:r>
        <w:t>Goa</w:t>
      </w:r>
      <w:r>
        <w:t>l</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="SpacingCharChar"/>
        <w:ind w:left="284"/>
        <w:rPr>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Orient in the documentation of the first four fundamental principles:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="SpacingCharChar"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:r>
        <w:t>Computer Language Programmed Within Itself</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
      </w:pPr>
      <w:r>
        <w:t>- Generic, Not Generated</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
      </w:pPr>
      <w:r>
        <w:t>- Small Code Base</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
      </w:pPr>
      <w:r>
        <w:t>- Everything Only (Lack Of Choice = Guarantees)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
      </w:pPr>
      <w:r>
        <w:t>Super-</w:t>
      </w:r>
      <w:r>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:t>roject</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>This is a sub-project of the super-project</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Circle Language Spec, Fundamental Principles Spec Part A</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
      </w:pPr>
      <w:r>
        <w:t>Date</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:r>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:t>ime</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
      </w:pPr>
      <w:r>
        <w:t>2008-03-05 – 2008-03-06</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="284"/>
      </w:pPr>
      <w:r>
        <w:t>2 hours</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> of work</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
      </w:pPr>
      <w:r>
        <w:t>Project Plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="478" w:hanging="194"/>
      </w:pPr>
      <w:r>
        <w:t>- Work on ‘Small Code Base’</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="478" w:hanging="194"/>
      </w:pPr>
      <w:r>
        <w:t>- Work on ‘Everything Only (Lack of Choice = Guarantees)’</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="478" w:hanging="194"/>
      </w:pPr>
      <w:r>
        <w:t>- Change order of four most basic principles</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="478" w:hanging="194"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">- Change name of </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Creator In Creator</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> principles to </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Computer Language Programmed Within Itself</w:t>
      </w:r>
    </w:p>
    <w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($bodyXml)

$s = $d.Styles("Heading3")
$pf = $s.ParagraphFormat
$pf.SpaceBefore = 9
$pf.SpaceAfter = 9
$f = $s.Font
$f.Size = 10
$f.Italic = $false

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
